# Applies the "Updated symbol list" commit: refreshed prices/volumes and
# re-ranked several coin rows (their Coin/Link/Price/Volume cells moved
# to a different row) on Sheet1 of cryptos.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell, new text, and whether it is a numeric-looking
# string (price / percentage) that must be forced to Text format first —
# otherwise Excel's auto-detection would coerce it into a real number
# (losing the exact displayed text / introducing float rounding).
$updates = @(
    @{ Cell = 'D2'; Value = '246.68'; AsText = $true }
    @{ Cell = 'E2'; Value = '0.97%'; AsText = $true }
    @{ Cell = 'D3'; Value = '30.33'; AsText = $true }
    @{ Cell = 'E3'; Value = '11.45%'; AsText = $true }
    @{ Cell = 'D4'; Value = '5.181'; AsText = $true }
    @{ Cell = 'E4'; Value = '0.40%'; AsText = $true }
    @{ Cell = 'D5'; Value = '0.05728'; AsText = $true }
    @{ Cell = 'E5'; Value = '1.66%'; AsText = $true }
    @{ Cell = 'B7'; Value = 'MXToken'; AsText = $false }
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; AsText = $false }
    @{ Cell = 'D7'; Value = '0.8559'; AsText = $true }
    @{ Cell = 'E7'; Value = '4.98%'; AsText = $true }
    @{ Cell = 'B8'; Value = 'FTXToken'; AsText = $false }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; AsText = $false }
    @{ Cell = 'D8'; Value = '0.8821'; AsText = $true }
    @{ Cell = 'E8'; Value = '5.97%'; AsText = $true }
    @{ Cell = 'B9'; Value = 'WazirX'; AsText = $false }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; AsText = $false }
    @{ Cell = 'D9'; Value = '0.1366'; AsText = $true }
    @{ Cell = 'E9'; Value = '2.65%'; AsText = $true }
    @{ Cell = 'B10'; Value = 'MandalaExchangeToken'; AsText = $false }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; AsText = $false }
    @{ Cell = 'D10'; Value = '0.07101'; AsText = $true }
    @{ Cell = 'E10'; Value = '2.75%'; AsText = $true }
    @{ Cell = 'B11'; Value = 'BitrueCoin'; AsText = $false }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; AsText = $false }
    @{ Cell = 'D11'; Value = '0.02870'; AsText = $true }
    @{ Cell = 'E11'; Value = '-2.15%'; AsText = $true }
    @{ Cell = 'B12'; Value = 'BitMartToken'; AsText = $false }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; AsText = $false }
    @{ Cell = 'D12'; Value = '0.09394'; AsText = $true }
    @{ Cell = 'E12'; Value = '-0.03%'; AsText = $true }
    @{ Cell = 'B13'; Value = 'BitForexToken'; AsText = $false }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; AsText = $false }
    @{ Cell = 'D13'; Value = '0.001523'; AsText = $true }
    @{ Cell = 'E13'; Value = '1.05%'; AsText = $true }
    @{ Cell = 'B14'; Value = 'CoinExToken'; AsText = $false }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'; AsText = $false }
    @{ Cell = 'D14'; Value = '0.04157'; AsText = $true }
    @{ Cell = 'E14'; Value = '-1.83%'; AsText = $true }
    @{ Cell = 'B15'; Value = 'One'; AsText = $false }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'; AsText = $false }
    @{ Cell = 'D15'; Value = '0.0006033'; AsText = $true }
    @{ Cell = 'E15'; Value = '0.42%'; AsText = $true }
    @{ Cell = 'B16'; Value = 'TigerCash'; AsText = $false }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'; AsText = $false }
    @{ Cell = 'D16'; Value = '0.005965'; AsText = $true }
    @{ Cell = 'E16'; Value = '-2.78%'; AsText = $true }
    @{ Cell = 'B17'; Value = 'LEO'; AsText = $false }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; AsText = $false }
    @{ Cell = 'D17'; Value = '3.490'; AsText = $true }
    @{ Cell = 'E17'; Value = '-1.22%'; AsText = $true }
    @{ Cell = 'B18'; Value = 'GateToken'; AsText = $false }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; AsText = $false }
    @{ Cell = 'D18'; Value = '3.073'; AsText = $true }
    @{ Cell = 'E18'; Value = '2.31%'; AsText = $true }
    @{ Cell = 'D19'; Value = '2.182'; AsText = $true }
    @{ Cell = 'E19'; Value = '-1.98%'; AsText = $true }
    @{ Cell = 'D20'; Value = '0.3223'; AsText = $true }
    @{ Cell = 'E20'; Value = '3.53%'; AsText = $true }
    @{ Cell = 'E21'; Value = '4.50%'; AsText = $true }
    @{ Cell = 'D22'; Value = '0.1301'; AsText = $true }
    @{ Cell = 'E22'; Value = '0.73%'; AsText = $true }
    @{ Cell = 'D23'; Value = '3.512'; AsText = $true }
    @{ Cell = 'E23'; Value = '-6.27%'; AsText = $true }
    @{ Cell = 'D24'; Value = '0.1381'; AsText = $true }
    @{ Cell = 'E24'; Value = '0.53%'; AsText = $true }
    @{ Cell = 'E25'; Value = '-1.02%'; AsText = $true }
    @{ Cell = 'D26'; Value = '0.004489'; AsText = $true }
    @{ Cell = 'E26'; Value = '0.30%'; AsText = $true }
    @{ Cell = 'E27'; Value = '23.48%'; AsText = $true }
    @{ Cell = 'D28'; Value = '0.0001384'; AsText = $true }
    @{ Cell = 'E28'; Value = '0.06%'; AsText = $true }
    @{ Cell = 'D40'; Value = '0.03782'; AsText = $true }
    @{ Cell = 'E40'; Value = '3.54%'; AsText = $true }
    @{ Cell = 'B41'; Value = 'BKEXToken'; AsText = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'; AsText = $false }
    @{ Cell = 'D41'; Value = '0.1072'; AsText = $true }
    @{ Cell = 'E41'; Value = '-21.96%'; AsText = $true }
    @{ Cell = 'B42'; Value = 'CEJI'; AsText = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'; AsText = $false }
    @{ Cell = 'D42'; Value = '0.002599'; AsText = $true }
    @{ Cell = 'E42'; Value = '-1.51%'; AsText = $true }
    @{ Cell = 'B43'; Value = 'KickToken'; AsText = $false }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'; AsText = $false }
    @{ Cell = 'D43'; Value = '0.003511'; AsText = $true }
    @{ Cell = 'E43'; Value = '-41.85%'; AsText = $true }
    @{ Cell = 'D44'; Value = '0.009998'; AsText = $true }
    @{ Cell = 'E44'; Value = '21.38%'; AsText = $true }
    @{ Cell = 'D45'; Value = '0.00005092'; AsText = $true }
    @{ Cell = 'E45'; Value = '-5.77%'; AsText = $true }
    @{ Cell = 'E46'; Value = '0.09%'; AsText = $true }
    @{ Cell = 'D47'; Value = '0.08004'; AsText = $true }
    @{ Cell = 'E47'; Value = '-40.69%'; AsText = $true }
    @{ Cell = 'D48'; Value = '0.002764'; AsText = $true }
    @{ Cell = 'E48'; Value = '4.31%'; AsText = $true }
    @{ Cell = 'D49'; Value = '0.00002101'; AsText = $true }
    @{ Cell = 'E49'; Value = '0.09%'; AsText = $true }
    @{ Cell = 'D50'; Value = '0.0002001'; AsText = $true }
    @{ Cell = 'E50'; Value = '0.09%'; AsText = $true }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
